# Updates betting-odds values on "Sheet1" per the FlashScore odds refresh
# (commit: "Atualizando o arquivo XLSX").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.3
$ws.Cells.Item(2, 9).Value = 3.7
$ws.Cells.Item(2, 12).Value = 4.75
$ws.Cells.Item(2, 14).Value = 4.75
$ws.Cells.Item(2, 15).Value = 1.8
$ws.Cells.Item(2, 16).Value = 1.91
$ws.Cells.Item(2, 17).Value = 3.6
$ws.Cells.Item(2, 18).Value = 1.29
$ws.Cells.Item(2, 19).Value = 1.8
$ws.Cells.Item(2, 20).Value = 2
$ws.Cells.Item(2, 24).Value = 9
$ws.Cells.Item(2, 25).Value = 12
$ws.Cells.Item(2, 26).Value = 23
$ws.Cells.Item(2, 30).Value = 6
$ws.Cells.Item(2, 34).Value = 15
$ws.Cells.Item(2, 35).Value = 15
$ws.Cells.Item(2, 41).Value = 17
$ws.Cells.Item(2, 50).Value = 26
$ws.Cells.Item(3, 7).Value = 1.95
$ws.Cells.Item(3, 8).Value = 3.1
$ws.Cells.Item(3, 9).Value = 4.75
$ws.Cells.Item(3, 10).Value = 2.75
$ws.Cells.Item(3, 12).Value = 5.5
$ws.Cells.Item(3, 13).Value = 1.13
$ws.Cells.Item(3, 14).Value = 6
$ws.Cells.Item(3, 21).Value = 2.5
$ws.Cells.Item(3, 22).Value = 1.5
$ws.Cells.Item(3, 23).Value = 4.75
$ws.Cells.Item(3, 24).Value = 7.5
$ws.Cells.Item(3, 26).Value = 17
$ws.Cells.Item(3, 27).Value = 21
$ws.Cells.Item(3, 30).Value = 6.5
$ws.Cells.Item(3, 31).Value = 23
$ws.Cells.Item(3, 33).Value = 8.5
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 35).Value = 19
$ws.Cells.Item(3, 36).Value = 51
$ws.Cells.Item(3, 37).Value = 51
$ws.Cells.Item(3, 38).Value = 67
$ws.Cells.Item(3, 40).Value = 3.6
$ws.Cells.Item(3, 41).Value = 12
$ws.Cells.Item(3, 42).Value = 29
$ws.Cells.Item(3, 43).Value = 41
$ws.Cells.Item(3, 48).Value = 101
$ws.Cells.Item(3, 49).Value = 6
$ws.Cells.Item(3, 50).Value = 29
$ws.Cells.Item(3, 52).Value = 126
$ws.Cells.Item(3, 53).Value = 201
$ws.Cells.Item(4, 7).Value = 1.7
$ws.Cells.Item(4, 8).Value = 3.3
$ws.Cells.Item(4, 9).Value = 6
$ws.Cells.Item(4, 10).Value = 2.4
$ws.Cells.Item(4, 12).Value = 6.5
$ws.Cells.Item(4, 15).Value = 1.5
$ws.Cells.Item(4, 16).Value = 2.5
$ws.Cells.Item(4, 17).Value = 2.6
$ws.Cells.Item(4, 18).Value = 1.48
$ws.Cells.Item(4, 26).Value = 12
$ws.Cells.Item(4, 30).Value = 7
$ws.Cells.Item(4, 33).Value = 11
$ws.Cells.Item(4, 34).Value = 29
$ws.Cells.Item(4, 41).Value = 9.5
$ws.Cells.Item(4, 42).Value = 26
$ws.Cells.Item(5, 12).Value = 3.1
$ws.Cells.Item(5, 21).Value = 1.95
$ws.Cells.Item(5, 22).Value = 1.8
$ws.Cells.Item(5, 26).Value = 34
$ws.Cells.Item(5, 34).Value = 10
$ws.Cells.Item(5, 41).Value = 19
$ws.Cells.Item(5, 42).Value = 29
$ws.Cells.Item(5, 49).Value = 4.33
$ws.Cells.Item(5, 53).Value = 81
$ws.Cells.Item(7, 24).Value = 9.5
$ws.Cells.Item(7, 52).Value = 67
$ws.Cells.Item(8, 7).Value = 2.7
$ws.Cells.Item(8, 12).Value = 3.75
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(10, 18).Value = 1.75
$ws.Cells.Item(10, 19).Value = 1.44
$ws.Cells.Item(10, 20).Value = 2.63
$ws.Cells.Item(10, 24).Value = 10
$ws.Cells.Item(10, 33).Value = 9.5
$ws.Cells.Item(10, 46).Value = 2.63
